$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "56.932.60"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +3.17%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.501.59"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.54%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.22%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "495.94"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +3.55%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.02"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +10.48%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.32%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.79%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.517.67"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.19%  "

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +6.35%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0995"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.73%  "

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +3.25%  "

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.30%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.936.80"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.39%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "57.031.02"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +3.13%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.51"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +5.07%  "

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.06%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.518.17"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.80%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.56"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +4.58%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.38"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +4.06%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "325.44"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +2.12%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.27%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.95"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +5.08%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "59.11"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +2.46%  "

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.75%  "

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.16%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.998"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.76%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.612.30"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.27%  "

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +4.70%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0817"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +4.98%  "

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.19%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "151.79"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.05%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.46"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.82%  "

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +4.21%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.27"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.95%  "

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +5.94%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.83"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +4.09%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.881"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +4.06%  "

$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.40"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +6.70%  "

$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "34.36"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.40%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0568"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +3.35%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.53"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +4.43%  "

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.29%  "

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.22%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.95"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +10.44%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "268.69"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +8.94%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0929"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +3.12%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0231"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +4.09%  "

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.69%  "

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +3.33%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.911.43"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -3.10%  "
